$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# --- Row 18: Clase Pedregal: atributos ---
$ws.Range("C18").Value = "Clase Pedregal: atributos"
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 0.00069444444444444447
$ws.Range("H18").Value = 0.55555555555555558023
$ws.Range("I18").Value = 0.55625000000000002220
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 5

# --- Row 19: Getter, Setters y Constructor ---
$ws.Range("C19").Value = "Getter, Setters y Constructor"
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 0.00208333333333333330
$ws.Range("H19").Value = 0.55625000000000002220
$ws.Range("I19").Value = 0.55833333333333334814
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 62

# --- Row 20: Método resolución ---
$ws.Range("C20").Value = "Método resolución"
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 0.00694444444444444406
$ws.Range("H20").Value = 0.55902777777777779011
$ws.Range("I20").Value = 0.56944444444444441977
$ws.Range("K20").Value = 4
$ws.Range("L20").Value = 0.00694444444444444406
$ws.Range("M20").Value = 35

# --- Row 21: Método casaCabeDesdeCasilla ---
$ws.Range("C21").Value = "Método casaCabeDesdeCasilla"
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 0.00208333333333333330
$ws.Range("H21").Value = 0.56944444444444441977
$ws.Range("I21").Value = 0.57083333333333330373
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.00347222222222222203
$ws.Range("M21").Value = 3

# --- Row 22: Método ubicarCasa ---
$ws.Range("C22").Value = "Método ubicarCasa"
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 0.00694444444444444406
$ws.Range("H22").Value = 0.57083333333333330373
$ws.Range("I22").Value = 0.57638888888888895057
$ws.Range("K22").Value = 4
$ws.Range("L22").Value = 0.00694444444444444406
$ws.Range("M22").Value = 10

# --- Row 23: Método orientaciónCasa ---
$ws.Range("C23").Value = "Método orientaciónCasa"
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 0.00694444444444444406
$ws.Range("H23").Value = 0.57638888888888895057
$ws.Range("I23").Value = 0.58124999999999993339
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 15

# --- Row 24: MétodoMostrarTerreno ---
$ws.Range("C24").Value = "MétodoMostrarTerreno"
$ws.Range("F24").Value = 10
$ws.Range("G24").Value = 0.00347222222222222203
$ws.Range("H24").Value = 0.58124999999999993339
$ws.Range("I24").Value = 0.58333333333333337034
$ws.Range("J24").Formula = '=IFERROR(IF(OR(ISBLANK(H24),ISBLANK(I24)),"",IF(I24>=H24,I24-H24,"Error")),"Error")'
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.00347222222222222203
$ws.Range("M24").Value = 11
$ws.Range("N24").Formula = '=IFERROR(IF(OR(J24="",ISBLANK(L24)),"",J24+L24),"Error")'

# --- Row 25: Main ---
$ws.Range("C25").Value = "Main"
$ws.Range("F25").Value = 10
$ws.Range("G25").Value = 0.00694444444444444406
$ws.Range("H25").Value = 0.58333333333333337034
$ws.Range("I25").Value = 0.58680555555555558023
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 6

# --- Row 26 totals label cleared ---
$ws.Range("B26").Value = ""

# --- Row 30: Diseño time tracking ---
$ws.Range("B30").Value = 0.01736111111111111188
$ws.Range("C30").Value = 0.59027777777777779011
$ws.Range("D30").Value = 0.60416666666666662966

# --- sheet view: selection moves to B33:D33 (and scroll back to A1) ---
$ws.Range("A1").Select()
$ws.Range("B33:D33").Select()

# Force a full recalculation so every dependent formula reflects final inputs.
$excel.CalculateFull()
